$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value even when the text looks numeric,
# without leaving a residual explicit NumberFormat on the cell.
function Set-TextValue($sheet, [string]$cellRef, [string]$val) {
    $rng = $sheet.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range('D2').Value = '26.202.79'
$ws.Range('E2').Value = '  -1.14%  '
$ws.Range('D3').Value = '1.658.40'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('E4').Value = '  +0.31%  '
Set-TextValue $ws 'D5' '217.08'
$ws.Range('E5').Value = '  -1.33%  '
Set-TextValue $ws 'D6' '0.5164'
$ws.Range('E6').Value = '  -2.15%  '
$ws.Range('E7').Value = '  +0.27%  '
Set-TextValue $ws 'D8' '0.2638'
$ws.Range('E8').Value = '  -1.49%  '
Set-TextValue $ws 'D9' '0.06264'
$ws.Range('E9').Value = '  -1.76%  '
Set-TextValue $ws 'D10' '20.74'
$ws.Range('E10').Value = '  -4.65%  '
Set-TextValue $ws 'D11' '0.07753'
$ws.Range('E11').Value = '  -0.70%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.673.17'
$ws.Range('E12').Value = '  +0.18%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws 'D13' '4.478'
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('D14').Value = '1.885.77'
$ws.Range('E14').Value = '  -0.78%  '
Set-TextValue $ws 'D15' '0.5452'
$ws.Range('E15').Value = '  -2.18%  '
$ws.Range('D16').Value = '0.0₅8121'
$ws.Range('E16').Value = '  -2.25%  '
Set-TextValue $ws 'D17' '64.80'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '26.209.46'
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('E19').Value = '  +0.34%  '
Set-TextValue $ws 'D20' '4.600'
$ws.Range('E20').Value = '  -3.48%  '
Set-TextValue $ws 'D21' '191.75'
$ws.Range('E21').Value = '  -0.76%  '
Set-TextValue $ws 'D22' '10.08'
$ws.Range('E22').Value = '  -2.30%  '
Set-TextValue $ws 'D23' '5.983'
$ws.Range('E23').Value = '  -5.17%  '
$ws.Range('E24').Value = '  +0.40%  '
Set-TextValue $ws 'D25' '139.47'
$ws.Range('E25').Value = '  +0.21%  '
Set-TextValue $ws 'D26' '0.1219'
Set-TextValue $ws 'D27' '7.263'
$ws.Range('E27').Value = '  -1.79%  '
Set-TextValue $ws 'D28' '16.11'
$ws.Range('E28').Value = '  -1.19%  '
Set-TextValue $ws 'D29' '1.441'
$ws.Range('E29').Value = '  +1.13%  '
Set-TextValue $ws 'D30' '0.05914'
$ws.Range('E30').Value = '  -4.72%  '
$ws.Range('E31').Value = '  -1.36%  '
Set-TextValue $ws 'D32' '3.542'
$ws.Range('E32').Value = '  -1.84%  '
$ws.Range('E33').Value = '  -4.51%  '
Set-TextValue $ws 'D34' '1.581'
$ws.Range('E34').Value = '  -6.18%  '
Set-TextValue $ws 'D35' '0.9584'
$ws.Range('E35').Value = '  -4.85%  '
$ws.Range('E36').Value = '  +0.29%  '
Set-TextValue $ws 'D37' '2.771'
$ws.Range('E37').Value = '  -0.25%  '
Set-TextValue $ws 'D38' '0.5659'
$ws.Range('E38').Value = '  -7.29%  '
Set-TextValue $ws 'D39' '6.047'
$ws.Range('E39').Value = '  +0.16%  '
Set-TextValue $ws 'D40' '0.01587'
$ws.Range('E40').Value = '  -1.76%  '
Set-TextValue $ws 'D41' '0.8521'
$ws.Range('E41').Value = '  -0.57%  '
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('D43').Value = '1.011.37'
$ws.Range('E43').Value = '  -7.02%  '
Set-TextValue $ws 'D44' '100.65'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').Value = '1.800.53'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('E46').Value = '  -3.35%  '
Set-TextValue $ws 'D47' '56.41'
$ws.Range('E47').Value = '  -3.37%  '
Set-TextValue $ws 'D48' '1.007'
$ws.Range('E48').Value = '  +0.33%  '
Set-TextValue $ws 'D49' '8.043'
$ws.Range('E49').Value = '  -0.65%  '
$ws.Range('E50').Value = '  -0.55%  '
Set-TextValue $ws 'D51' '0.4221'
$ws.Range('E51').Value = '  -0.17%  '
